$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 17 (ID 16)
$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Patient"
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 8.3
$ws.Range("E17").Value = 4.2
$ws.Range("F17").Value = 8.1
$ws.Range("G17").Value = 16.3
$ws.Range("H17").Value = 8.2
$ws.Range("I17").Value = 3.2
$ws.Range("J17").Value = 16.5
$ws.Range("K17").Value = 0.1
$ws.Range("L17").Value = 3.6
$ws.Range("M17").Value = 0.1
$ws.Range("N17").Value = 8.3
$ws.Range("O17").Value = 0.2
$ws.Range("P17").Value = 0.1
$ws.Range("Q17").Value = 0.1
$ws.Range("R17").Value = 16.6
$ws.Range("S17").Value = 16.6
$ws.Range("T17").Value = 5.5
$ws.Range("U17").Value = 1
$ws.Range("V17").Value = "JT"
$ws.Range("W17").Value = 2
$ws.Range("X17").Value = "Should do research in ICU Should follow up all the patients Should see how patients are getting on"

# New row 18 (ID 17)
$ws.Range("A18").Value = 17
$ws.Range("B18").Value = "Patient"
$ws.Range("D18").Value = 0.1
$ws.Range("E18").Value = 7.8
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 16.5
$ws.Range("I18").Value = 16.4
$ws.Range("J18").Value = 16.5
$ws.Range("K18").Value = 0.3
$ws.Range("L18").Value = 1.1
$ws.Range("M18").Value = 0.4
$ws.Range("N18").Value = 8.2
$ws.Range("O18").Value = 15.7
$ws.Range("P18").Value = 0.5
$ws.Range("Q18").Value = 15
$ws.Range("R18").Value = 15
$ws.Range("S18").Value = 16.2
$ws.Range("T18").Value = 15
$ws.Range("U18").Value = 1
$ws.Range("V18").Value = "JT"
$ws.Range("W18").Value = 4

# New row 19 (ID 18)
$ws.Range("A19").Value = 18
$ws.Range("B19").Value = "Patient"
$ws.Range("D19").Value = 16.4
$ws.Range("E19").Value = 16.6
$ws.Range("F19").Value = 16.6
$ws.Range("G19").Value = 16.6
$ws.Range("H19").Value = 16.6
$ws.Range("I19").Value = 12.6
$ws.Range("J19").Value = 16.3
$ws.Range("K19").Value = 16.6
$ws.Range("L19").Value = 0.2
$ws.Range("M19").Value = 7.9
$ws.Range("N19").Value = 3.2
$ws.Range("O19").Value = 0.2
$ws.Range("P19").Value = 8
$ws.Range("Q19").Value = 8.2
$ws.Range("R19").Value = 0.1
$ws.Range("S19").Value = 12.2
$ws.Range("T19").Value = 0.2
$ws.Range("U19").Value = 1
$ws.Range("V19").Value = "JT"
$ws.Range("W19").Value = 4
$ws.Range("X19").Value = "If there’s anything we can do we should do it"

$ws.Range("X19").Select()
